# LMDI: Versjon 1.0.6 6e05e801b0fc67a31d9121f33125496b6f7ed95a
# - clear the "Experimental" value ("true") in B7
# - bump the "Date" value in B8 from 2025-03-10 to 2025-09-12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 = Experimental / true  -> clear the value cell (B7) while keeping formatting
$ws.Range("B7").Value = ""

# Row 8 = Date / 2025-03-10 -> 2025-09-12
# Force text (not an Excel date serial) by pre-formatting as Text, then
# restore the original "General" formatting via a formats-only paste from a
# neighboring cell that already carries the worksheet's standard style.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2025-09-12"
$ws.Range("A8").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
